$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# generators.py now inserts a blank row (row 4) between the previous
# day's orders and the orders appended on the current day, then the
# new order is appended starting at row 5.

# DATE and PHONE look numeric/date-like to Excel's auto-detection, so
# mark them as Text before writing so they are stored as real text
# (matching the source data) instead of being auto-converted into a
# date serial number / plain number.
$ws.Range("A5").NumberFormat = "@"
$ws.Range("D5").NumberFormat = "@"

$ws.Range("A5").Value = "01/11/2024"
$ws.Range("B5").Value = "SO240111001"
$ws.Range("C5").Value = "Ashley Baker"
$ws.Range("D5").Value = "9157994875"
$ws.Range("E5").Value = "Smashing Pumpkins"
$ws.Range("F5").Value = "Siamese Dream"
$ws.Range("G5").Value = 55
$ws.Range("H5").Value = 120.99
$ws.Range("I5").Value = "AEC"
$ws.Range("J5").Value = "DVD"
$ws.Range("K5").Value = "Ashley"
